# Update Name of Algo
# Applies updated numeric values to the result_data_RandomForest worksheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.14089999999999
$ws.Range("E6").Value = 12.39909999999999
$ws.Range("E7").Value = 11.99629999999999
$ws.Range("C8").Value = -11.9145
$ws.Range("E8").Value = 13.50370000000001
$ws.Range("B12").Value = 5.8079
$ws.Range("C12").Value = -14.56050000000002
$ws.Range("C14").Value = -11.8092
$ws.Range("E19").Value = 13.0896
$ws.Range("E21").Value = 12.70879999999999
$ws.Range("C22").Value = -10.96729999999999
$ws.Range("E24").Value = 12.70189999999998
